# Append a new row (row 10) with date 2025-03-11 to each of the 9 price
# sheets in the workbook, carrying forward the same price as the prior
# day (row 9), matching the "Updated Argent prices" commit.

$wb = $excel.ActiveWorkbook

$newDate = "2025-03-11"

# Sheet name -> price value to place in column B of the new row.
$sheetPrices = [ordered]@{
    "N-Dense"                   = "40"
    "N-Type"                    = "43"
    "N-type Wafer"               = "1.19"
    "Cell Topcon 183mm"          = "0.295"
    "Module Topcon 183mm"        = "0.1"
    "Silver Rear_side"           = "5,273"
    "Silver Busbar front-side"   = "7,895"
    "Silver finger front-side"   = "7,945"
    "USD_CNY"                    = "7.2597"
}

foreach ($sheetName in $sheetPrices.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find the first empty row right after the existing data (row 10 for
    # every sheet here, since each currently holds rows 1-9).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $newRow = $lastRow + 1

    # Use a leading apostrophe so Excel stores these as plain text
    # (matching the existing inline-string cells) instead of converting
    # the date-looking / numeric-looking text into a real date or number.
    $ws.Cells.Item($newRow, 1).Value = "'" + $newDate
    $ws.Cells.Item($newRow, 1).Style = "Normal"

    $ws.Cells.Item($newRow, 2).Value = "'" + $sheetPrices[$sheetName]
    $ws.Cells.Item($newRow, 2).Style = "Normal"
}
